# Appends the cumulative module status matrix block (21 new paragraphs:
# a blank line, a "---" rule, a title, an "Updated:" line, a blank line,
# a tab-delimited header row, and 15 tab-delimited module status rows)
# immediately after the existing final paragraph, before the section break.
# Word's Range.InsertXML (WordprocessingML "single file package" form) is
# used so the exact run formatting (Helvetica Light, 12pt/24 half-points)
# and `xml:space="preserve"` text runs are preserved verbatim.

$d = $word.ActiveDocument

# A zero-length (collapsed) range at the very end of the body content -
# i.e. immediately after the last existing paragraph's paragraph mark and
# before <w:sectPr/>. Inserting XML here appends new paragraphs without
# disturbing any existing content.
$insertionPoint = $d.Range($d.Content.End, $d.Content.End)

$newContentPackageXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">---</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Cumulative Main Requirement Status Matrix (Do Not Remove Previous Entries)</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Updated: 2026-02-18</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve"/></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Module Name</w:t><w:tab/><w:t xml:space="preserve">Developed</w:t><w:tab/><w:t xml:space="preserve">Partial Developed</w:t><w:tab/><w:t xml:space="preserve">Need To Develop</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">User / Merchant Management</w:t><w:tab/><w:t xml:space="preserve">Merchant onboarding flow (register-&gt;OTP placeholder-&gt;plan-&gt;store setup), tenant resolver middleware, core RBAC scaffold, team invite/accept, lifecycle APIs/UI, franchise/backoffice base models</w:t><w:tab/><w:t xml:space="preserve">Platform owner seed/delegation flows, custom role builder UX, lifecycle bulk ops and approval chains</w:t><w:tab/><w:t xml:space="preserve">Business onboarding automation, full franchise hierarchy workflows, enterprise role governance</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Storefront System (Retail + Wholesale)</w:t><w:tab/><w:t xml:space="preserve">Landing + auth/onboarding screens, responsive frontend base</w:t><w:tab/><w:t xml:space="preserve">Theme/layout pages exist but not full storefront CMS, B2B visibility toggles not complete</w:t><w:tab/><w:t xml:space="preserve">Theme marketplace, drag-drop homepage builder, navigation builder, static page CMS, hybrid retail/wholesale runtime controls</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Product Management</w:t><w:tab/><w:t xml:space="preserve">Products/categories/variants/media core entities + CRUD APIs wired</w:t><w:tab/><w:t xml:space="preserve">Bulk upload, advanced inventory rules, wholesale MOQ/tier pricing partial placeholders</w:t><w:tab/><w:t xml:space="preserve">Complete wholesale rule engine (MOQ/pack/case/tier), Excel import/export, richer media/video pipeline</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Order Management</w:t><w:tab/><w:t xml:space="preserve">Orders/order-items core entities + CRUD APIs, lifecycle status fields, notes basics</w:t><w:tab/><w:t xml:space="preserve">Invoice/GST flow partial, refund/cancel flow partial, shipment integrations partial</w:t><w:tab/><w:t xml:space="preserve">Complete unified order engine (retail/wholesale/social/manual), GST invoice generation, refunds + courier integrations + audit-grade transitions</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Customer Management</w:t><w:tab/><w:t xml:space="preserve">Customers + addresses + CRUD APIs wired</w:t><w:tab/><w:t xml:space="preserve">Customer groups/segmentation and B2B credit ledger partial</w:t><w:tab/><w:t xml:space="preserve">Guest checkout model, GSTIN enforcement rules, customer-specific pricing matrix, full segmentation workflows</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Payments &amp; Financial Features</w:t><w:tab/><w:t xml:space="preserve">Payment plugin interface scaffold, billing plan/subscription models scaffold</w:t><w:tab/><w:t xml:space="preserve">Real gateways, tokenized card/wallet storage, UPI/COD flows incomplete</w:t><w:tab/><w:t xml:space="preserve">Production gateway plugins, partial payment/refund settlement flows, reconciliation, payouts, GST finance reporting</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Store Builder &amp; Customization</w:t><w:tab/><w:t xml:space="preserve">Branding/auth UI customization done for onboarding/auth pages</w:t><w:tab/><w:t xml:space="preserve">Theme controls minimal and not merchant self-serve</w:t><w:tab/><w:t xml:space="preserve">Live preview editor, section builder, custom CSS/JS sandbox, banner/slider/store identity management</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Plugin / App Ecosystem</w:t><w:tab/><w:t xml:space="preserve">Plugin architecture direction started (payment plugin registry)</w:t><w:tab/><w:t xml:space="preserve">No install/uninstall marketplace flow yet</w:t><w:tab/><w:t xml:space="preserve">Full app marketplace, permission scopes, plugin billing, developer portal, integration categories</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Logistics &amp; Fulfillment</w:t><w:tab/><w:t xml:space="preserve">Data model hooks available via existing order domain</w:t><w:tab/><w:t xml:space="preserve">No full logistics module yet</w:t><w:tab/><w:t xml:space="preserve">Shipping rule engine, courier APIs, pincode serviceability, labels, returns, wholesale freight workflows</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Marketing &amp; Growth</w:t><w:tab/><w:t xml:space="preserve">Basic platform shell only</w:t><w:tab/><w:t xml:space="preserve">Discount/coupon foundations not fully implemented</w:t><w:tab/><w:t xml:space="preserve">Discount engine, flash sales, combo offers, abandoned cart recovery, email/SMS/WhatsApp automation, SEO/social tools</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Reporting &amp; Analytics</w:t><w:tab/><w:t xml:space="preserve">Audit logs viewer + key admin logs</w:t><w:tab/><w:t xml:space="preserve">Sales/traffic/product analytics dashboards partial</w:t><w:tab/><w:t xml:space="preserve">Full BI dashboards, GST/tax reports, conversion funnels, exportable analytics</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Security &amp; Compliance</w:t><w:tab/><w:t xml:space="preserve">Opaque tokens, Turnstile verify wiring, rate-limit hooks, tenancy isolation, policy-based auth, audit logs</w:t><w:tab/><w:t xml:space="preserve">MFA/WebAuthn partial config only, CSRF hardening partial, DataProtection key persistence pending</w:t><w:tab/><w:t xml:space="preserve">Production-grade MFA/WebAuthn rollout, KMS encryption for secrets, backup/recovery controls, fraud scoring + incident workflows</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Advanced / Competitive Features</w:t><w:tab/><w:t xml:space="preserve">Architecture allows future extensibility</w:t><w:tab/><w:t xml:space="preserve">No AI/multi-channel modules yet</w:t><w:tab/><w:t xml:space="preserve">AI assistants, smart pricing/fraud/catalog, WhatsApp/POS/marketplace channel connectors</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">SaaS Business Engine</w:t><w:tab/><w:t xml:space="preserve">Billing plan/subscription base models + some admin endpoints</w:t><w:tab/><w:t xml:space="preserve">Usage metering and feature gating partial</w:t><w:tab/><w:t xml:space="preserve">Trial automation, usage limits, upgrade/downgrade, add-on monetization, subscription ops dashboard</w:t></w:r></w:p>
<w:p><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">Technical Architecture</w:t><w:tab/><w:t xml:space="preserve">Backend .NET + PostgreSQL + React frontend wired; Docker/Render deployment files present</w:t><w:tab/><w:t xml:space="preserve">Read/write split, replication strategy, infra automation partial</w:t><w:tab/><w:t xml:space="preserve">Production cloud architecture (AWS/Azure), HA, observability, backups, scaling playbooks</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$insertionPoint.InsertXML($newContentPackageXml)

Write-Output ("Paragraphs after insert: " + $d.Paragraphs.Count)
